$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H43").Value = 1038.8572
$ws.Range("I43").Value = 1043.75
$ws.Range("J43").Value = 1032.3334
$ws.Range("K43").Value = 1043.75
$ws.Range("L43").Value = 1032.3334
$ws.Range("M43").Value = -974.75
$ws.Range("N43").Value = -1170.3334

$ws.Range("H88").Value = 2400
$ws.Range("J88").Value = 2400
$ws.Range("L88").Value = 2400
$ws.Range("N88").Value = -3212

$ws.Range("H91").Value = 2400
$ws.Range("J91").Value = 2400
$ws.Range("L91").Value = 2400
$ws.Range("N91").Value = -5208

$ws.Range("H116").Value = 6760.375
$ws.Range("I116").Value = 3749.5
$ws.Range("K116").Value = 3749.5
$ws.Range("M116").Value = -307.5

$ws.Range("H125").Value = 700.5
$ws.Range("I125").Value = 700.5
$ws.Range("K125").Value = 6304.5
$ws.Range("M125").Value = -3844.5

$ws.Range("H129").Value = 3616.818
$ws.Range("J129").Value = 3174.125
$ws.Range("L129").Value = 9522.375
$ws.Range("N129").Value = -19522.375

$ws.Range("H132").Value = 3059.5715
$ws.Range("I132").Value = 3543.6
$ws.Range("J132").Value = 1849.5
$ws.Range("K132").Value = 10630.8
$ws.Range("L132").Value = 5548.5
$ws.Range("M132").Value = -8100.799999999999
$ws.Range("N132").Value = -10608.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2041.85
$ws.Range("I2").Value = 1553.3636
$ws.Range("J2").Value = 2638.889
$ws.Range("K2").Value = 1553.3636
$ws.Range("L2").Value = 2638.889
$ws.Range("M2").Value = -1440.3636
$ws.Range("N2").Value = -2864.889

$ws.Range("H32").Value = 3083033.5
$ws.Range("I32").Value = 3185174.8
$ws.Range("J32").Value = 2333997.8
$ws.Range("K32").Value = 3185174.8
$ws.Range("L32").Value = 2333997.8
$ws.Range("M32").Value = -3184887.8
$ws.Range("N32").Value = -2334571.8

$ws.Range("H74").Value = 2085.1667
$ws.Range("I74").Value = 1837
$ws.Range("K74").Value = 1837
$ws.Range("M74").Value = -963

$ws.Range("H77").Value = 2085.1667
$ws.Range("I77").Value = 1837
$ws.Range("K77").Value = 9185
$ws.Range("M77").Value = -4817

$ws.Range("H116").Value = 2041.85
$ws.Range("I116").Value = 1553.3636
$ws.Range("J116").Value = 2638.889
$ws.Range("K116").Value = 1553.3636
$ws.Range("L116").Value = 2638.889
$ws.Range("M116").Value = 740.6364000000001
$ws.Range("N116").Value = -7226.889

$ws.Range("H132").Value = 1241

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2041.85
$ws.Range("I3").Value = 1553.3636
$ws.Range("J3").Value = 2638.889
$ws.Range("K3").Value = 1553.3636
$ws.Range("L3").Value = 2638.889
$ws.Range("M3").Value = -1439.3636
$ws.Range("N3").Value = -2866.889

$ws.Range("H92").Value = 47000
$ws.Range("J92").Value = 47000
$ws.Range("L92").Value = 47000
$ws.Range("N92").Value = -51992

$ws.Range("H99").Value = 1995
$ws.Range("I99").Value = 2005
$ws.Range("J99").Value = 1988.3334
$ws.Range("K99").Value = 2005
$ws.Range("L99").Value = 1988.3334
$ws.Range("M99").Value = -507
$ws.Range("N99").Value = -4984.3334

$ws.Range("H107").Value = 1666.4
$ws.Range("I107").Value = 1296.1111
$ws.Range("J107").Value = 4999
$ws.Range("K107").Value = 1296.1111
$ws.Range("L107").Value = 4999
$ws.Range("M107").Value = 623.8888999999999
$ws.Range("N107").Value = -8839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3280.8
$ws.Range("I16").Value = 3634.6667
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 3634.6667
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -3347.6667
$ws.Range("N16").Value = -3324

$ws.Range("H31").Value = 2082.4285
$ws.Range("I31").Value = 2071.75
$ws.Range("J31").Value = 2096.6667
$ws.Range("K31").Value = 2071.75
$ws.Range("L31").Value = 2096.6667
$ws.Range("M31").Value = -1776.75
$ws.Range("N31").Value = -2686.6667

$ws.Range("H34").Value = 2082.4285
$ws.Range("I34").Value = 2071.75
$ws.Range("J34").Value = 2096.6667
$ws.Range("K34").Value = 2071.75
$ws.Range("L34").Value = 2096.6667
$ws.Range("M34").Value = -1869.75
$ws.Range("N34").Value = -2500.6667

$ws.Range("H86").Value = 9223.25
$ws.Range("I86").Value = 5632.3335
$ws.Range("K86").Value = 5632.3335
$ws.Range("M86").Value = -4509.3335

$ws.Range("H89").Value = 9223.25
$ws.Range("I89").Value = 5632.3335
$ws.Range("K89").Value = 28161.6675
$ws.Range("M89").Value = -22545.6675

$ws.Range("H94").Value = 142604.5
$ws.Range("I94").Value = 186601.83
$ws.Range("J94").Value = 10612.5
$ws.Range("K94").Value = 186601.83
$ws.Range("L94").Value = 10612.5
$ws.Range("M94").Value = -186150.83
$ws.Range("N94").Value = -11514.5

$ws.Range("H105").Value = 2904.9167
$ws.Range("I105").Value = 2329.3333
$ws.Range("K105").Value = 2329.3333
$ws.Range("M105").Value = -582.3332999999998

$ws.Range("H107").Value = 922.875
$ws.Range("I107").Value = 813.8333
$ws.Range("K107").Value = 813.8333
$ws.Range("M107").Value = 1106.1667

$ws.Range("H113").Value = 3280.8
$ws.Range("I113").Value = 3634.6667
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 3634.6667
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = -1464.6667
$ws.Range("N113").Value = -7090

$ws.Range("H122").Value = 2060.8333
$ws.Range("I122").Value = 2266.25
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 6798.75
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -4348.75
$ws.Range("N122").Value = -9850

$ws.Range("H134").Value = 2442.9473
$ws.Range("I134").Value = 2377.4119
$ws.Range("K134").Value = 7132.2357
$ws.Range("M134").Value = -4597.2357

$ws.Range("H141").Value = 240200
$ws.Range("J141").Value = 240200
$ws.Range("L141").Value = 240200
$ws.Range("N141").Value = -250560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H6").Value = 278.5
$ws.Range("I6").Value = 284.2
$ws.Range("K6").Value = 852.5999999999999
$ws.Range("M6").Value = -739.5999999999999

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H132").Value = 51499.5
$ws.Range("I132").Value = 51499.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 463495.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -460965.5
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 34999.4
$ws.Range("J15").Value = 34999.4
$ws.Range("L15").Value = 34999.4
$ws.Range("N15").Value = -35575.4

$ws.Range("H17").Value = 46.25
$ws.Range("I17").Value = 46.25
$ws.Range("K17").Value = 46.25
$ws.Range("M17").Value = 121.75

$ws.Range("H81").Value = 34999.4
$ws.Range("J81").Value = 34999.4
$ws.Range("L81").Value = 34999.4
$ws.Range("N81").Value = -36995.4

$ws.Range("H84").Value = 34999.4
$ws.Range("J84").Value = 34999.4
$ws.Range("L84").Value = 104998.2
$ws.Range("N84").Value = -114982.2

$ws.Range("H102").Value = 2347
$ws.Range("I102").Value = 2396.5715
$ws.Range("K102").Value = 2396.5715
$ws.Range("M102").Value = -774.5715

$ws.Range("H132").Value = 4937.75
$ws.Range("I132").Value = 5372.2354
$ws.Range("J132").Value = 2475.6667
$ws.Range("K132").Value = 16116.7062
$ws.Range("L132").Value = 7427.000100000001
$ws.Range("M132").Value = -13586.7062
$ws.Range("N132").Value = -12487.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1924.6666
$ws.Range("I22").Value = 2560
$ws.Range("J22").Value = 1470.8572
$ws.Range("K22").Value = 2560
$ws.Range("L22").Value = 1470.8572
$ws.Range("M22").Value = -2265
$ws.Range("N22").Value = -2060.8572

$ws.Range("H27").Value = 1924.6666
$ws.Range("I27").Value = 2560
$ws.Range("J27").Value = 1470.8572
$ws.Range("K27").Value = 2560
$ws.Range("L27").Value = 1470.8572
$ws.Range("M27").Value = -2453
$ws.Range("N27").Value = -1684.8572

$ws.Range("H40").Value = 3789.1333
$ws.Range("I40").Value = 3504.1
$ws.Range("K40").Value = 3504.1
$ws.Range("M40").Value = -3368.1

$ws.Range("H55").Value = 1348.8235
$ws.Range("I55").Value = 1258.375
$ws.Range("K55").Value = 1258.375
$ws.Range("M55").Value = -1085.375

$ws.Range("H100").Value = 4252.933
$ws.Range("J100").Value = 2798.6667
$ws.Range("L100").Value = 2798.6667
$ws.Range("N100").Value = -3880.6667

$ws.Range("H132").Value = 3487.875
$ws.Range("J132").Value = 3520.8
$ws.Range("L132").Value = 10562.4
$ws.Range("N132").Value = -15622.4

$ws.Range("H136").Value = 1250
$ws.Range("I136").Value = 1250
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3750
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1200
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1896.8334
$ws.Range("I126").Value = 1395.75
$ws.Range("K126").Value = 4187.25
$ws.Range("M126").Value = -1717.25

$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
